$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2007年" row (row 2) was removed from the source data, shifting
# every subsequent row (2010/2012/2015/2017年) up by one and shrinking
# the sheet's used range from A1:Y6 to A1:Y5.
$ws.Rows.Item(2).Delete()
